# 2.1.1.1e.xlsx update: add a new "2020" data column (N) and refresh the
# latest-year (2019, column M / L9 / L12) figures with revised numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Revise a handful of already-existing 2019 figures (column M, plus
#    two 2018 corrections in L9/L12) to match the republished data.
# ---------------------------------------------------------------------
$ws.Range("M5").Value  = 68.4
$ws.Range("M6").Value  = 108.2
$ws.Range("M7").Value  = 51.7
$ws.Range("M8").Value  = 97.7
$ws.Range("L9").Value  = 105.6
$ws.Range("M9").Value  = 106.7
$ws.Range("M10").Value = 124.2
$ws.Range("M11").Value = 138.8
$ws.Range("L12").Value = 27.1
$ws.Range("M12").Value = 33.9
$ws.Range("M13").Value = 96
$ws.Range("M14").Value = 7.7

# ---------------------------------------------------------------------
# 2. Add the new column N (year 2020). Copy the formatting from the
#    matching cell one column to the left (M) so the new cells reuse the
#    same styles (number format / font / borders) instead of creating
#    new style entries. Row 3 has no "M"-like counterpart value wise, so
#    its border-only style is copied from A14, which carries the same
#    bottom-border/no-numfmt style used for that separator row.
# ---------------------------------------------------------------------
$ws.Range("A14").Copy()
$ws.Range("N3").PasteSpecial(-4122)

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)

$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)

$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)

$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)

$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)

$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial(-4122)

$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial(-4122)

$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial(-4122)

$ws.Range("M13").Copy()
$ws.Range("N13").PasteSpecial(-4122)

$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fill in the 2020 values for the new column.
# ---------------------------------------------------------------------
$ws.Range("N4").Value  = 2020
$ws.Range("N5").Value  = 68.5
$ws.Range("N6").Value  = 106.7
$ws.Range("N7").Value  = 53.2
$ws.Range("N8").Value  = 49.6
$ws.Range("N9").Value  = 108.9
$ws.Range("N10").Value = 107.8
$ws.Range("N11").Value = 155.7
$ws.Range("N12").Value = 25.9
$ws.Range("N13").Value = 103.5
$ws.Range("N14").Value = 11

# ---------------------------------------------------------------------
# 4. Set the page setup (print) properties that were added for this
#    sheet.
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
